# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# The workbook's Sheet1 holds per-game pitching stats for Germán Márquez (2021).
# Column G (header "K") previously held the total pitch count that resulted in a
# strike ("Strike#"); it is being regenerated to hold the actual strikeout (K)
# totals for each outing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new strikeout (K) value, taken from the regenerated
# save_data (row 2 = most recent game, row 38 = oldest game in the sheet).
$kValues = @{
    2  = 4
    3  = 4
    4  = 1
    5  = 7
    6  = 5
    7  = 2
    8  = 3
    9  = 4
    10 = 6
    11 = 7
    12 = 9
    13 = 10
    14 = 5
    15 = 9
    16 = 11
    17 = 5
    18 = 7
    19 = 2
    20 = 5
    21 = 6
    22 = 7
    23 = 6
    24 = 8
    25 = 8
    26 = 6
    27 = 1
    28 = 5
    29 = 8
    30 = 6
    31 = 7
    32 = 5
    33 = 2
    34 = 1
    35 = 4
    36 = 6
    37 = 3
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
